# Jordan Burdett gets his username appended, the same way the other
# authors already have theirs ("Leanne Capewell - leanneec",
# "Benedict Hobart - bhobart"). Word's "_GoBack" bookmark (the marker
# for "last edit point", used by Shift+F5) moves along with the edit -
# it needs to be removed from its old spot (after the "merging of DNA"
# paragraph) and re-created around the newly-typed text.

$d = $word.ActiveDocument

# --- 1. Find the "Jordan Burdett -" byline paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    # Range.Text carries a trailing paragraph-mark char (chr 13) - strip
    # it before comparing against the plain byline text.
    $candidateText = $candidate.Range.Text.TrimEnd([char]13, [char]7)
    if ($candidateText -eq "Jordan Burdett -") {
        $target = $candidate
        break
    }
}

if ($target -ne $null) {

    # --- 2. _GoBack is a singleton bookmark: Word relocates it to
    #        wherever you last typed, so the copy sitting after the
    #        "merging of DNA" paragraph has to go before a new one is
    #        created around the text we're about to insert below.
    try {
        $goBack = $d.Bookmarks("_GoBack")
        $goBack.Delete()
    } catch {
        # no pre-existing _GoBack bookmark - nothing to clean up
    }

    # --- 3. Rebuild the paragraph with " jburdett" appended: a plain
    #        (non-bold) run for the space, the username wrapped in
    #        spell-check proofing marks (it isn't a dictionary word,
    #        same as "leanneec"/"bhobart" elsewhere in this doc), and
    #        the relocated _GoBack bookmark spanning that new text,
    #        exactly like Word leaves behind after you type it
    #        interactively.
    $r = $target.Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' +
        '<w:p w14:paraId="46AF9E92" w14:textId="77777777" w:rsidR="00C27DC6" w:rsidRDefault="000D2379">' +
        '<w:pPr><w:pStyle w:val="normal0"/><w:ind w:right="-89"/><w:jc w:val="center"/></w:pPr>' +
        '<w:r><w:rPr><w:b/></w:rPr><w:t>Jordan Burdett -</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>jburdett</w:t></w:r>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '</w:p>' +
        '</w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $null = $r.InsertXML($xml)
}
